$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Ltf"
$ws.Range("C2").Value = "Lrp11"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01868033333333333
$ws.Range("H2").Value = 0.056041
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.172479
$ws.Range("N2").Value = 0.5174369999999999
$ws.Range("O2").Value = 0.1112535186193337
$ws.Range("P2").Value = 0.1112535186193337
$ws.Range("Q2").Value = 0.003221965212999999
$ws.Range("R2").Value = 0.02899768691699999
$ws.Range("S2").Value = 0.1112535186193337
$ws.Range("T2").Value = 0.1112535186193337

$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Ltf"
$ws.Range("C3").Value = "Lrp11"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01868033333333333
$ws.Range("H3").Value = 0.056041
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.9644846666666668
$ws.Range("N3").Value = 2.893454
$ws.Range("O3").Value = 0.6221181292856633
$ws.Range("P3").Value = 0.6221181292856633
$ws.Range("Q3").Value = 0.01801689506822222
$ws.Range("R3").Value = 0.162152055614
$ws.Range("S3").Value = 0.6221181292856633
$ws.Range("T3").Value = 0.6221181292856633

$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Ltf"
$ws.Range("C4").Value = "Lrp11"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01868033333333333
$ws.Range("H4").Value = 0.056041
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4133603333333333
$ws.Range("N4").Value = 1.240081
$ws.Range("O4").Value = 0.2666283520950029
$ws.Range("P4").Value = 0.2666283520950029
$ws.Range("Q4").Value = 0.007721708813444444
$ws.Range("R4").Value = 0.069495379321
$ws.Range("S4").Value = 0.2666283520950029
$ws.Range("T4").Value = 0.2666283520950029
